$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Thang Duong Chi" -> "Thang" (spell-check-flagged) + " Duong Chi"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Thang Duong Chi", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null

$xml1 = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>
<w:p w:rsidR='00C64291' w:rsidRDefault='001275F9' w:rsidP='00B5208A'>
  <w:pPr>
    <w:jc w:val='center'/>
  </w:pPr>
  <w:proofErr w:type='spellStart'/>
  <w:r><w:t>Thang</w:t></w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r><w:t xml:space='preserve'> Duong Chi</w:t></w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$rng1.InsertXML($xml1) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: append " = " (separate run) after the plain (non-bold) "Server
# Data" that introduces the "Data Structure Traveling with the Flow" value.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Server Data", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null
$rng2.Find.Execute("Server Data", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "", 0) | Out-Null

$para2 = $rng2.Paragraphs(1)
$pRng2 = $para2.Range

$xml2 = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>
<w:p w:rsidR='004A3EA4' w:rsidRPr='004A3EA4' w:rsidRDefault='001275F9' w:rsidP='004A3EA4'>
  <w:pPr>
    <w:rPr><w:lang w:val='en-US'/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Server Data</w:t></w:r>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> = </w:t></w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$pRng2.InsertXML($xml2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: add a centered tab stop + trailing tab to the "volume/time"
# comment paragraph, then append a brand-new "Server Data data structure"
# paragraph (carrying the _GoBack bookmark that used to close the old last
# paragraph).
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("The volume/time may not be correct.", $true, $false, `
                    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$para3 = $rng3.Paragraphs(1)
$pRng3 = $para3.Range
# Trim off the trailing paragraph mark so the two replacement paragraphs
# don't leave a stray empty paragraph behind.
$pRng3b = $d.Range($pRng3.Start, $pRng3.End - 1)

$xml3 = @"
<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>
<w:p w:rsidR='004A3EA4' w:rsidRPr='004A3EA4' w:rsidRDefault='001275F9' w:rsidP='004A3EA4'>
  <w:pPr>
    <w:tabs><w:tab w:val='center' w:pos='4513'/></w:tabs>
    <w:rPr><w:lang w:val='en-US'/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>The volume/time may not be correct.</w:t></w:r>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:tab/></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:tabs><w:tab w:val='center' w:pos='4513'/></w:tabs>
    <w:rPr><w:lang w:val='en-US'/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'>Server Data </w:t></w:r>
  <w:proofErr w:type='spellStart'/>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>data</w:t></w:r>
  <w:proofErr w:type='spellEnd'/>
  <w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t xml:space='preserve'> structure</w:t></w:r>
  <w:bookmarkStart w:id='0' w:name='_GoBack'/>
  <w:bookmarkEnd w:id='0'/>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$pRng3b.InsertXML($xml3) | Out-Null
